$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet 1 updates ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 06:37:24"
$ws1.Cells.Item(3, 1).Value = "Total filas: 68"
$ws1.Cells.Item(46, 1).Value = "06:37:24"
$ws1.Cells.Item(46, 2).Value = "07:01"
$ws1.Cells.Item(46, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(46, 4).Value = 24
$ws1.Cells.Item(47, 1).Value = "05:47:32"
$ws1.Cells.Item(47, 2).Value = "07:04"
$ws1.Cells.Item(47, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(47, 4).Value = 77
$ws1.Cells.Item(48, 1).Value = "05:18:23"
$ws1.Cells.Item(48, 3).Value = "15_ABASTO"
$ws1.Cells.Item(48, 4).Value = 107
$ws1.Cells.Item(49, 1).Value = "06:02:16"
$ws1.Cells.Item(49, 2).Value = "07:05"
$ws1.Cells.Item(49, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(49, 4).Value = 63
$ws1.Cells.Item(50, 2).Value = "07:07"
$ws1.Cells.Item(50, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(50, 4).Value = 109
$ws1.Cells.Item(51, 2).Value = "07:11"
$ws1.Cells.Item(51, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(51, 4).Value = 113
$ws1.Cells.Item(52, 1).Value = "05:18:23"
$ws1.Cells.Item(52, 2).Value = "07:15"
$ws1.Cells.Item(52, 4).Value = 117
$ws1.Cells.Item(53, 1).Value = "06:37:24"
$ws1.Cells.Item(53, 2).Value = "07:16"
$ws1.Cells.Item(53, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(53, 4).Value = 39
$ws1.Cells.Item(54, 2).Value = "07:16"
$ws1.Cells.Item(54, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(54, 4).Value = 74
$ws1.Cells.Item(55, 2).Value = "07:21"
$ws1.Cells.Item(55, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(55, 4).Value = 94
$ws1.Cells.Item(56, 1).Value = "06:02:16"
$ws1.Cells.Item(56, 2).Value = "07:23"
$ws1.Cells.Item(56, 3).Value = "10_OLMOS"
$ws1.Cells.Item(56, 4).Value = 81
$ws1.Cells.Item(57, 2).Value = "07:27"
$ws1.Cells.Item(57, 3).Value = "10_OLMOS"
$ws1.Cells.Item(57, 4).Value = 100
$ws1.Cells.Item(58, 2).Value = "07:31"
$ws1.Cells.Item(58, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(58, 4).Value = 104
$ws1.Cells.Item(59, 1).Value = "05:47:32"
$ws1.Cells.Item(59, 2).Value = "07:31"
$ws1.Cells.Item(59, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(59, 4).Value = 104
$ws1.Cells.Item(60, 1).Value = "06:02:16"
$ws1.Cells.Item(60, 2).Value = "07:32"
$ws1.Cells.Item(60, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(60, 4).Value = 90
$ws1.Cells.Item(61, 1).Value = "05:47:32"
$ws1.Cells.Item(61, 2).Value = "07:32"
$ws1.Cells.Item(61, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(61, 4).Value = 105
$ws1.Cells.Item(62, 1).Value = "06:37:24"
$ws1.Cells.Item(62, 2).Value = "07:34"
$ws1.Cells.Item(62, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(62, 4).Value = 57
$ws1.Cells.Item(63, 1).Value = "05:47:32"
$ws1.Cells.Item(63, 2).Value = "07:36"
$ws1.Cells.Item(63, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(63, 4).Value = 109
$ws1.Cells.Item(64, 2).Value = "07:37"
$ws1.Cells.Item(64, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(64, 4).Value = 95
$ws1.Cells.Item(65, 1).Value = "05:47:32"
$ws1.Cells.Item(65, 2).Value = "07:39"
$ws1.Cells.Item(65, 3).Value = "10_OLMOS"
$ws1.Cells.Item(65, 4).Value = 112
$ws1.Cells.Item(65, 5).Value = "LP1912"
$ws1.Cells.Item(66, 1).Value = "06:37:24"
$ws1.Cells.Item(66, 2).Value = "07:47"
$ws1.Cells.Item(66, 3).Value = "14_ABASTO"
$ws1.Cells.Item(66, 4).Value = 70
$ws1.Cells.Item(66, 5).Value = "LP1912"
$ws1.Cells.Item(67, 1).Value = "06:02:16"
$ws1.Cells.Item(67, 2).Value = "07:48"
$ws1.Cells.Item(67, 3).Value = "14_ABASTO"
$ws1.Cells.Item(67, 4).Value = 106
$ws1.Cells.Item(67, 5).Value = "LP1912"
$ws1.Cells.Item(68, 1).Value = "06:02:16"
$ws1.Cells.Item(68, 2).Value = "07:51"
$ws1.Cells.Item(68, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(68, 4).Value = 109
$ws1.Cells.Item(68, 5).Value = "LP1912"
$ws1.Cells.Item(69, 1).Value = "06:37:24"
$ws1.Cells.Item(69, 2).Value = "08:12"
$ws1.Cells.Item(69, 3).Value = "15_ABASTO"
$ws1.Cells.Item(69, 4).Value = 95
$ws1.Cells.Item(69, 5).Value = "LP1912"
$ws1.Cells.Item(70, 1).Value = "06:37:24"
$ws1.Cells.Item(70, 2).Value = "08:21"
$ws1.Cells.Item(70, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(70, 4).Value = 104
$ws1.Cells.Item(70, 5).Value = "LP1912"
$ws1.Cells.Item(71, 1).Value = "06:37:24"
$ws1.Cells.Item(71, 2).Value = "08:22"
$ws1.Cells.Item(71, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(71, 4).Value = 105
$ws1.Cells.Item(71, 5).Value = "LP1912"
$ws1.Cells.Item(72, 1).Value = "06:37:24"
$ws1.Cells.Item(72, 2).Value = "08:23"
$ws1.Cells.Item(72, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(72, 4).Value = 106
$ws1.Cells.Item(72, 5).Value = "LP1912"
$ws1.Cells.Item(73, 1).Value = "06:37:24"
$ws1.Cells.Item(73, 2).Value = "08:27"
$ws1.Cells.Item(73, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(73, 4).Value = 110
$ws1.Cells.Item(73, 5).Value = "LP1912"

# --- Sheet 2 updates ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 06:37:24"
$ws2.Cells.Item(3, 1).Value = "Total filas: 13"
$ws2.Cells.Item(18, 1).Value = "06:37:24"
$ws2.Cells.Item(18, 2).Value = "08:23"
$ws2.Cells.Item(18, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(18, 4).Value = 106
$ws2.Cells.Item(18, 5).Value = "LP1912"

# --- Sheet 3 updates ---
$ws3.Cells.Item(2, 1).Value = "Última actualización: 06:37:24"
$ws3.Cells.Item(3, 1).Value = "Total filas: 11"
$ws3.Cells.Item(15, 1).Value = "06:37:24"
$ws3.Cells.Item(15, 2).Value = "08:07"
$ws3.Cells.Item(15, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(15, 4).Value = 90
$ws3.Cells.Item(15, 5).Value = "L6203"
$ws3.Cells.Item(16, 1).Value = "06:37:24"
$ws3.Cells.Item(16, 2).Value = "08:30"
$ws3.Cells.Item(16, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(16, 4).Value = 113
$ws3.Cells.Item(16, 5).Value = "L6173"
